$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Arkusz1")

# Fix difficulty level (poz_trud, column C) values
$ws.Range("C6:C11").Value = 1
$ws.Range("C18:C23").Value = 2

# Fix correct-answer flag (praw, column G) values
$ws.Range("G20:G23").Value = "F"

# Adjust column widths (D -> 47.140625 chars, E -> 59.5703125 chars)
$ws.Columns.Item(4).ColumnWidth = 46.25
$ws.Columns.Item(5).ColumnWidth = 58.6

# Update view / selection state
$ws.Range("B18").Select()
